$d = $word.ActiveDocument

# 1. Brian's paragraph: append a second sentence about strategy.
$d.Content.Find.Execute(
    "Brian – Completed the full implementation for the Enemy’s movement method",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Brian – Completed the full implementation for the Enemy’s movement method and helped devise strategy for implementation.",
    2)

# 2. Jasim's paragraph: append details about meetings and design logic.
$d.Content.Find.Execute(
    "Jasim – Coded the game timer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Jasim – Coded the game timer, scheduled group meetings, and assisted in game design logic.",
    2)

# 3. Vera's paragraph: "Started" -> "Spearheaded" and add "production" before the comma.
$d.Content.Find.Execute(
    "Started the report and code, giving a good backbone for the others to assist.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Spearheaded the report and code production, giving a good backbone for the others to assist.",
    2)

# 4. Challenges paragraph: "external libraries" -> "core libraries".
$d.Content.Find.Execute(
    "getting used to and using the external libraries",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "getting used to and using the core libraries",
    2)
